{"js": "// \"NTT Address\" block on the first page used to read, across two\n// centered paragraphs:\n//   Cluj Napoca\n//   City, State 400158\n// Merge them into a single centered paragraph:\n//   Cluj Napoca City, 400158\n// (the first paragraph's text is replaced and the second paragraph,\n// now redundant, is removed).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet cityParagraph = null;\nlet addressParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (text === \"Cluj Napoca\") {\n    const next = paragraphs.items[i + 1];\n    if (next.text.trim() === \"City, State 400158\") {\n      cityParagraph = para;\n      addressParagraph = next;\n      break;\n    }\n  }\n}\n\nif (cityParagraph && addressParagraph) {\n  cityParagraph.insertText(\"Cluj Napoca City, 400158\", \"Replace\");\n  addressParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# \"NTT Address\" block on the first page used to read, across two\n# centered paragraphs:\n#   Cluj Napoca\n#   City, State 400158\n# Merge them into a single centered paragraph:\n#   Cluj Napoca City, 400158\n# (the first paragraph's text is replaced and the second paragraph,\n# now redundant, is removed).\n\n$d = $word.ActiveDocument\n\n$cityIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -lt $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($text -eq \"Cluj Napoca\") {\n        $nextText = $d.Paragraphs.Item($i + 1).Range.Text.Trim()\n        if ($nextText -eq \"City, State 400158\") {\n            $cityIndex = $i\n            break\n        }\n    }\n}\n\nif ($cityIndex -ge 1) {\n    $cityParagraph = $d.Paragraphs.Item($cityIndex)\n    $cityParagraph.Range.Text = \"Cluj Napoca City, 400158\"\n\n    $addressParagraph = $d.Paragraphs.Item($cityIndex + 1)\n    $addressParagraph.Range.Delete()\n}\n"}
